$d = $word.ActiveDocument

# --- First paragraph: "**ID__AFFARS_5325_topic_20__ID** " -----------------

# 1. Replace the placeholder ID text. Matching the original text together
#    with its trailing space (which lived in a second, separate run) lets
#    Word's Find/Replace consume both runs and leave a single run behind
#    containing just the new ID text, with no extra trailing space.
$null = $d.Content.Find.Execute("**ID__AFFARS_5325_topic_20__ID** ", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "**ID__AFFARS_5325_7501__ID**", 2)

# 2. Give the paragraph a (line-less) paragraph border on all four sides,
#    reserving 5pt of space on each edge, and increase its left indent
#    from 6pt (120 twips) to 11.25pt (225 twips).
$p1 = $d.Paragraphs.Item(1)
$b1 = $p1.Format.Borders
$b1.DistanceFromTop = 5
$b1.DistanceFromLeft = 5
$b1.DistanceFromBottom = 5
$b1.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25
